$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.972.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.302.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.99%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.31%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("E11").Value = "  -0.66%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.46%  "

$ws.Range("E13").Value = "  +0.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.973"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.15%  "

$ws.Range("E15").Value = "  -4.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.645.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.295.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.046.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.31%  "

$ws.Range("E20").Value = "  -1.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.74%  "

$ws.Range("E22").Value = "  -6.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.96%  "

$ws.Range("E26").Value = "  +0.42%  "

$ws.Range("E27").Value = "  -3.88%  "

$ws.Range("E28").Value = "  +3.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.90%  "

$ws.Range("E31").Value = "  +1.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0897"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("E33").Value = "  -5.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.16%  "

$ws.Range("E35").Value = "  +10.70%  "

$ws.Range("E36").Value = "  -1.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.46%  "

$ws.Range("E38").Value = "  -0.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.68%  "

$ws.Range("E40").Value = "  -3.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.86%  "

$ws.Range("E43").Value = "  -2.57%  "

$ws.Range("E44").Value = "  -3.95%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.29%  "

$ws.Range("E48").Value = "  -1.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.568.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.35%  "
